# Apply updated "dSF" (and one "dS0") values as part of a data repull / recalculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new value for column F (dSF), except row 31 which updates column E (dS0)
$ws.Range("F3").Value = -5
$ws.Range("F9").Value = 4
$ws.Range("F10").Value = -1
$ws.Range("F11").Value = -1
$ws.Range("F17").Value = 0
$ws.Range("F20").Value = 3
$ws.Range("F22").Value = 1
$ws.Range("F27").Value = -1
$ws.Range("F28").Value = -1
$ws.Range("E31").Value = 0
$ws.Range("F35").Value = -1
$ws.Range("F43").Value = -4
$ws.Range("F45").Value = -1
$ws.Range("F46").Value = 2
$ws.Range("F48").Value = 1
$ws.Range("F49").Value = 1
$ws.Range("F53").Value = 0
$ws.Range("F54").Value = -2
$ws.Range("F55").Value = 5
